$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276206556/TJaAJ2AJpb-camera-video-segment-driver-1743276209056.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCGWPQSKN%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T170017Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIHrT8Ax277UsZsbA0IqxVgGgU8AyUelzUphbpv6M7DxYAiEA%2FhjKOtrool3%2BkLphEkVV9Vergmv0QwpQadb2yFZRMycq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFs%2BOwIcf%2B7rUuDHLSq6A1%2FI6LOkoQCMSJo7IUCMq9UfQrz3zXqwkOQoeV20VDLWPxWn3gtzBXDBO4gdClkkVJY%2B88vXPtIdxE%2BRlrA%2B9xNat97EpgaHdKqr%2FwnsC2OP76ESt5C1bltzBS7%2BbRTkU21TpM%2BOD3lHMH%2BizC7HwVfpfGvj%2B5UFm7UZ5%2Bc7mb7FKW9SWPXjj7OqIBsLQS5OCD1fsZxJaH%2Bib5tbMPzdN29kL%2F7mfH2yso28g2HDHvpmSocwKKTvCNRBdNTxbFOPVss6JqEognXrRaP1Xi1PG4LdaB8aLsuGe2xBZpdY0i3ldjY%2B%2F75LF%2FrXZ188SlrjPpS%2FTQ%2B%2BO5WPFKpd%2Fh4hEBEPCHeEOa5V9EGoErTn3Ff732c8htk9K6wfveNRDDJ%2B7ej4DWjcmjTJCnQYLOdR9p5uHCsMKSpFND1WtxeVtUFIElELxXkIDTzJp1Xoy7MMEdC7CPTmNsuolfSnjfp4U9KHlggL%2FpG4CzYKguANcyzxyCEw18PTrdMo3UNKe5SjgzIEjv4LUQVnzwHANzJ8Kpn8tBe%2BRi%2BJprdBye9M0HqijVCITzioI%2B5Bv9S4rkU9EkzpERGfuW%2FI55IwkNSlvwY6pQEgFjwnHYGnunyE6saWRGCl7BiC1EGCLjUdEaSD7Yp2wpwNx9QYlxqKwfZWyCJi4HK9MLlrn5ejDG%2BVTJETD92YD3%2BjWNeh92Q5c7Jw4ccuZI8F8B10wO4got34C0r3QMn0NdGcEA3Qc416WjjtTfYNPQTO0jsoztIWfqVxhoJ9fOJtqdbn6BQP6DGtKBQYK5onpz%2FkjXN8HKiIcxDeu2zrYvjqqpI%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2001%3A00%3A17%20GMT&X-Amz-Signature=3d145a8a33fa20939613b5a2a80fbad7f5219c28d5594952157404d5004b917e'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276164863/OF6jgx1Lcp-camera-video-segment-driver-1743276167363.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCGWPQSKN%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T170017Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIHrT8Ax277UsZsbA0IqxVgGgU8AyUelzUphbpv6M7DxYAiEA%2FhjKOtrool3%2BkLphEkVV9Vergmv0QwpQadb2yFZRMycq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFs%2BOwIcf%2B7rUuDHLSq6A1%2FI6LOkoQCMSJo7IUCMq9UfQrz3zXqwkOQoeV20VDLWPxWn3gtzBXDBO4gdClkkVJY%2B88vXPtIdxE%2BRlrA%2B9xNat97EpgaHdKqr%2FwnsC2OP76ESt5C1bltzBS7%2BbRTkU21TpM%2BOD3lHMH%2BizC7HwVfpfGvj%2B5UFm7UZ5%2Bc7mb7FKW9SWPXjj7OqIBsLQS5OCD1fsZxJaH%2Bib5tbMPzdN29kL%2F7mfH2yso28g2HDHvpmSocwKKTvCNRBdNTxbFOPVss6JqEognXrRaP1Xi1PG4LdaB8aLsuGe2xBZpdY0i3ldjY%2B%2F75LF%2FrXZ188SlrjPpS%2FTQ%2B%2BO5WPFKpd%2Fh4hEBEPCHeEOa5V9EGoErTn3Ff732c8htk9K6wfveNRDDJ%2B7ej4DWjcmjTJCnQYLOdR9p5uHCsMKSpFND1WtxeVtUFIElELxXkIDTzJp1Xoy7MMEdC7CPTmNsuolfSnjfp4U9KHlggL%2FpG4CzYKguANcyzxyCEw18PTrdMo3UNKe5SjgzIEjv4LUQVnzwHANzJ8Kpn8tBe%2BRi%2BJprdBye9M0HqijVCITzioI%2B5Bv9S4rkU9EkzpERGfuW%2FI55IwkNSlvwY6pQEgFjwnHYGnunyE6saWRGCl7BiC1EGCLjUdEaSD7Yp2wpwNx9QYlxqKwfZWyCJi4HK9MLlrn5ejDG%2BVTJETD92YD3%2BjWNeh92Q5c7Jw4ccuZI8F8B10wO4got34C0r3QMn0NdGcEA3Qc416WjjtTfYNPQTO0jsoztIWfqVxhoJ9fOJtqdbn6BQP6DGtKBQYK5onpz%2FkjXN8HKiIcxDeu2zrYvjqqpI%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2001%3A00%3A17%20GMT&X-Amz-Signature=a434aa8c3b06c11e5b8c9046f45e51f4991d4ff4952bd0c1b4fdee62c8c74637'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743276108738/gAxSKLP636-camera-video-segment-driver-1743276111238.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCGWPQSKN%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T170017Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIHrT8Ax277UsZsbA0IqxVgGgU8AyUelzUphbpv6M7DxYAiEA%2FhjKOtrool3%2BkLphEkVV9Vergmv0QwpQadb2yFZRMycq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFs%2BOwIcf%2B7rUuDHLSq6A1%2FI6LOkoQCMSJo7IUCMq9UfQrz3zXqwkOQoeV20VDLWPxWn3gtzBXDBO4gdClkkVJY%2B88vXPtIdxE%2BRlrA%2B9xNat97EpgaHdKqr%2FwnsC2OP76ESt5C1bltzBS7%2BbRTkU21TpM%2BOD3lHMH%2BizC7HwVfpfGvj%2B5UFm7UZ5%2Bc7mb7FKW9SWPXjj7OqIBsLQS5OCD1fsZxJaH%2Bib5tbMPzdN29kL%2F7mfH2yso28g2HDHvpmSocwKKTvCNRBdNTxbFOPVss6JqEognXrRaP1Xi1PG4LdaB8aLsuGe2xBZpdY0i3ldjY%2B%2F75LF%2FrXZ188SlrjPpS%2FTQ%2B%2BO5WPFKpd%2Fh4hEBEPCHeEOa5V9EGoErTn3Ff732c8htk9K6wfveNRDDJ%2B7ej4DWjcmjTJCnQYLOdR9p5uHCsMKSpFND1WtxeVtUFIElELxXkIDTzJp1Xoy7MMEdC7CPTmNsuolfSnjfp4U9KHlggL%2FpG4CzYKguANcyzxyCEw18PTrdMo3UNKe5SjgzIEjv4LUQVnzwHANzJ8Kpn8tBe%2BRi%2BJprdBye9M0HqijVCITzioI%2B5Bv9S4rkU9EkzpERGfuW%2FI55IwkNSlvwY6pQEgFjwnHYGnunyE6saWRGCl7BiC1EGCLjUdEaSD7Yp2wpwNx9QYlxqKwfZWyCJi4HK9MLlrn5ejDG%2BVTJETD92YD3%2BjWNeh92Q5c7Jw4ccuZI8F8B10wO4got34C0r3QMn0NdGcEA3Qc416WjjtTfYNPQTO0jsoztIWfqVxhoJ9fOJtqdbn6BQP6DGtKBQYK5onpz%2FkjXN8HKiIcxDeu2zrYvjqqpI%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2001%3A00%3A17%20GMT&X-Amz-Signature=36324cd4a2bfaaa117c540a78c23e7f1fac7f0e472fa26fe15d0063221c84d27'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205341/1743265652740/hq4516Mzws-camera-video-segment-driver-1743265655240.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCGWPQSKN%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T170017Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIHrT8Ax277UsZsbA0IqxVgGgU8AyUelzUphbpv6M7DxYAiEA%2FhjKOtrool3%2BkLphEkVV9Vergmv0QwpQadb2yFZRMycq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFs%2BOwIcf%2B7rUuDHLSq6A1%2FI6LOkoQCMSJo7IUCMq9UfQrz3zXqwkOQoeV20VDLWPxWn3gtzBXDBO4gdClkkVJY%2B88vXPtIdxE%2BRlrA%2B9xNat97EpgaHdKqr%2FwnsC2OP76ESt5C1bltzBS7%2BbRTkU21TpM%2BOD3lHMH%2BizC7HwVfpfGvj%2B5UFm7UZ5%2Bc7mb7FKW9SWPXjj7OqIBsLQS5OCD1fsZxJaH%2Bib5tbMPzdN29kL%2F7mfH2yso28g2HDHvpmSocwKKTvCNRBdNTxbFOPVss6JqEognXrRaP1Xi1PG4LdaB8aLsuGe2xBZpdY0i3ldjY%2B%2F75LF%2FrXZ188SlrjPpS%2FTQ%2B%2BO5WPFKpd%2Fh4hEBEPCHeEOa5V9EGoErTn3Ff732c8htk9K6wfveNRDDJ%2B7ej4DWjcmjTJCnQYLOdR9p5uHCsMKSpFND1WtxeVtUFIElELxXkIDTzJp1Xoy7MMEdC7CPTmNsuolfSnjfp4U9KHlggL%2FpG4CzYKguANcyzxyCEw18PTrdMo3UNKe5SjgzIEjv4LUQVnzwHANzJ8Kpn8tBe%2BRi%2BJprdBye9M0HqijVCITzioI%2B5Bv9S4rkU9EkzpERGfuW%2FI55IwkNSlvwY6pQEgFjwnHYGnunyE6saWRGCl7BiC1EGCLjUdEaSD7Yp2wpwNx9QYlxqKwfZWyCJi4HK9MLlrn5ejDG%2BVTJETD92YD3%2BjWNeh92Q5c7Jw4ccuZI8F8B10wO4got34C0r3QMn0NdGcEA3Qc416WjjtTfYNPQTO0jsoztIWfqVxhoJ9fOJtqdbn6BQP6DGtKBQYK5onpz%2FkjXN8HKiIcxDeu2zrYvjqqpI%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2001%3A00%3A17%20GMT&X-Amz-Signature=79bf1b1564a91a9453f48ef92696e52766967782b5c71ef7cb70fe292bca719d'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743255999969/Yjq7ZuVo4p-camera-video-segment-driver-1743256002469.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCGWPQSKN%2F20250330%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250330T170017Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECgaCXVzLXdlc3QtMiJHMEUCIHrT8Ax277UsZsbA0IqxVgGgU8AyUelzUphbpv6M7DxYAiEA%2FhjKOtrool3%2BkLphEkVV9Vergmv0QwpQadb2yFZRMycq5gMIkf%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDFs%2BOwIcf%2B7rUuDHLSq6A1%2FI6LOkoQCMSJo7IUCMq9UfQrz3zXqwkOQoeV20VDLWPxWn3gtzBXDBO4gdClkkVJY%2B88vXPtIdxE%2BRlrA%2B9xNat97EpgaHdKqr%2FwnsC2OP76ESt5C1bltzBS7%2BbRTkU21TpM%2BOD3lHMH%2BizC7HwVfpfGvj%2B5UFm7UZ5%2Bc7mb7FKW9SWPXjj7OqIBsLQS5OCD1fsZxJaH%2Bib5tbMPzdN29kL%2F7mfH2yso28g2HDHvpmSocwKKTvCNRBdNTxbFOPVss6JqEognXrRaP1Xi1PG4LdaB8aLsuGe2xBZpdY0i3ldjY%2B%2F75LF%2FrXZ188SlrjPpS%2FTQ%2B%2BO5WPFKpd%2Fh4hEBEPCHeEOa5V9EGoErTn3Ff732c8htk9K6wfveNRDDJ%2B7ej4DWjcmjTJCnQYLOdR9p5uHCsMKSpFND1WtxeVtUFIElELxXkIDTzJp1Xoy7MMEdC7CPTmNsuolfSnjfp4U9KHlggL%2FpG4CzYKguANcyzxyCEw18PTrdMo3UNKe5SjgzIEjv4LUQVnzwHANzJ8Kpn8tBe%2BRi%2BJprdBye9M0HqijVCITzioI%2B5Bv9S4rkU9EkzpERGfuW%2FI55IwkNSlvwY6pQEgFjwnHYGnunyE6saWRGCl7BiC1EGCLjUdEaSD7Yp2wpwNx9QYlxqKwfZWyCJi4HK9MLlrn5ejDG%2BVTJETD92YD3%2BjWNeh92Q5c7Jw4ccuZI8F8B10wO4got34C0r3QMn0NdGcEA3Qc416WjjtTfYNPQTO0jsoztIWfqVxhoJ9fOJtqdbn6BQP6DGtKBQYK5onpz%2FkjXN8HKiIcxDeu2zrYvjqqpI%3D&X-Amz-SignedHeaders=host&response-expires=Mon%2C%2031%20Mar%202025%2001%3A00%3A17%20GMT&X-Amz-Signature=0698f0138f92939d66121a2c63c7e07c9cc428a018f764d26af23afd0b5a7385'
